$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F3').Value = 680
$ws.Range('F9').Value = 7426
$ws.Range('F11').Value = 153
$ws.Range('F19').Value = 675
$ws.Range('F21').Value = 1171
$ws.Range('F23').Value = 566
$ws.Range('F24').Value = 10
$ws.Range('F26').Value = 36
$ws.Range('F27').Value = 564
$ws.Range('F29').Value = 4779
$ws.Range('G29').Value = 70
$ws.Range('F30').Value = 2297
$ws.Range('F31').Value = 3976
$ws.Range('F32').Value = 2216
$ws.Range('F37').Value = 56
$ws.Range('F42').Value = 541
$ws.Range('F43').Value = 284
$ws.Range('F45').Value = 853
$ws.Range('F48').Value = 21

# Sheet: 演出
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F4').Value = 15
$ws.Range('F15').Value = 361
$ws.Range('F20').Value = 69
$ws.Range('F25').Value = 110
$ws.Range('F32').Value = 1622
$ws.Range('F33').Value = 1622
$ws.Range('F35').Value = 23

# Sheet: 本地生活
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F6').Value = 1873
$ws.Range('F7').Value = 1889
$ws.Range('F10').Value = 1166
$ws.Range('F11').Value = 38
$ws.Range('F12').Value = 467
$ws.Range('F13').Value = 1881
$ws.Range('F14').Value = 8305
$ws.Range('F15').Value = 584

# Sheet: 全部类型
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 680
$ws.Range('F5').Value = 1873
$ws.Range('F8').Value = 7426
$ws.Range('F10').Value = 1166
$ws.Range('F11').Value = 15
$ws.Range('F12').Value = 467
$ws.Range('F13').Value = 1881
$ws.Range('F20').Value = 675
$ws.Range('F22').Value = 1171
$ws.Range('F25').Value = 361
$ws.Range('F28').Value = 566
$ws.Range('F29').Value = 10
$ws.Range('F31').Value = 69
$ws.Range('F44').Value = 110
$ws.Range('F46').Value = 541
$ws.Range('F47').Value = 285
$ws.Range('F51').Value = 1622

# Sheet: 全部类型 - row shift block (new event inserted, LookLook entry removed)
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('B33').Value = '2024-09-15'
$ws.Range('C33').Value = '【会员购提前抢】上海·宫村优子粉丝见面会'
$ws.Range('D33').Value = '西藏南路1号 上海大世界'
$ws.Range('E33').Value = '2024.09.15 12:00-09.16 17:30'
$ws.Range('F33').Value = 36
$ws.Range('G33').Value = 198
$ws.Range('H33').Value = 'https://show.bilibili.com/platform/detail.html?id=91139'
$ws.Range('I33').Value = '//i0.hdslb.com/bfs/openplatform/202408/i3owWPKf1724302904937.jpeg'

$ws.Range('B34').Value = '2024-09-15'
$ws.Range('C34').Value = '上海·GH·第五人格同人ONLY 1.0'
$ws.Range('D34').Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Range('E34').Value = '2024.09.15 10:00-09.15 17:00'
$ws.Range('F34').Value = 84
$ws.Range('G34').Value = 68
$ws.Range('H34').Value = 'https://show.bilibili.com/platform/detail.html?id=90638'
$ws.Range('I34').Value = '//i0.hdslb.com/bfs/openplatform/202408/plVDxJKi1723102207272.jpeg'

$ws.Range('B35').Value = '2024-09-15'
$ws.Range('C35').Value = '上海·SCGE动漫游戏嘉年华'
$ws.Range('D35').Value = '军工路1076号 纪希片场(秀场)'
$ws.Range('E35').Value = '2024.09.15 10:00-09.16 17:00'
$ws.Range('F35').Value = 4779
$ws.Range('G35').Value = 70
$ws.Range('H35').Value = 'https://show.bilibili.com/platform/detail.html?id=89993'
$ws.Range('I35').Value = '//i0.hdslb.com/bfs/openplatform/202408/aIJyQziE1723434354531.jpeg'

$ws.Range('B36').Value = '2024-09-15'
$ws.Range('C36').Value = '上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华'
$ws.Range('D36').Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range('E36').Value = '2024.09.15 11:00-09.16 16:00'
$ws.Range('F36').Value = 2297
$ws.Range('G36').Value = 65.8
$ws.Range('H36').Value = 'https://show.bilibili.com/platform/detail.html?id=90990'
$ws.Range('I36').Value = '//i1.hdslb.com/bfs/openplatform/202408/DutuUgvA1724127081751.jpeg'

$ws.Range('B37').Value = '2024-09-15'
$ws.Range('C37').Value = '上海·原神ONLY逐月节·原神&崩铁&崩三&绝区零·同人动漫嘉年华'
$ws.Range('D37').Value = '杨树浦路198号(金茂北外滩)B1层 Terra Park北外滩'
$ws.Range('E37').Value = '2024.09.15 09:30-09.17 22:30'
$ws.Range('F37').Value = 3976
$ws.Range('G37').Value = 78
$ws.Range('H37').Value = 'https://show.bilibili.com/platform/detail.html?id=89712'
$ws.Range('I37').Value = '//i2.hdslb.com/bfs/openplatform/202407/e9g9lWiy1721904672057.jpeg'

$ws.Range('B38').Value = '2024-09-15'
$ws.Range('C38').Value = '上海·城市动漫节2th'
$ws.Range('D38').Value = '西藏南路1号 上海大世界'
$ws.Range('E38').Value = '2024.09.15 10:00-09.16 18:00'
$ws.Range('F38').Value = 2216
$ws.Range('G38').Value = 68
$ws.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=89186'
$ws.Range('I38').Value = '//i1.hdslb.com/bfs/openplatform/202408/dyznHqyF1723780926438.jpeg'

$ws.Range('B39').Value = '2024-09-15'
$ws.Range('C39').Value = '上海·无限流同人only'
$ws.Range('D39').Value = '呼青路158号 交运智慧湾科创园25号楼'
$ws.Range('E39').Value = '2024.09.15 10:00-09.15 17:00'
$ws.Range('F39').Value = 232
$ws.Range('G39').Value = 89
$ws.Range('H39').Value = 'https://show.bilibili.com/platform/detail.html?id=90108'
$ws.Range('I39').Value = '//i2.hdslb.com/bfs/openplatform/202407/P3XVrcMn1722407440627.jpeg'

$ws.Range('B40').Value = '2024-09-15'
$ws.Range('C40').Value = '上海·第二届妖妖动漫游戏展'
$ws.Range('D40').Value = '吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙'
$ws.Range('E40').Value = '2024.09.15 10:00-09.16 17:00'
$ws.Range('F40').Value = 1159
$ws.Range('G40').Value = 68
$ws.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=90284'
$ws.Range('I40').Value = '//i2.hdslb.com/bfs/openplatform/202408/Q3xelO9p1722578696753.jpeg'

Write-Host "edits applied"